$d = $word.ActiveDocument

# Locate the "СПИСОК ИСПОЛНИТЕЛЕЙ" paragraph via Find (as text can repeat,
# Find gives us the matching Range directly).
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$found = $find.Execute("СПИСОК ИСПОЛНИТЕЛЕЙ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Expand to the whole paragraph (including the paragraph mark) so the
    # bold formatting is recorded on both the run and the paragraph mark's
    # run properties, matching a "Bold" toolbar toggle applied to the
    # paragraph.
    $para = $rng.Paragraphs(1).Range
    $para.Font.Bold = $true
    $para.Font.BoldBi = $true
}
